$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 0.5712403333333333
$ws.Range("N2").Value = 1.713721
$ws.Range("O2").Value = 0.1938341213320902
$ws.Range("P2").Value = 0.1963983345368606
$ws.Range("Q2").Value = 0.03770338530755556
$ws.Range("R2").Value = 0.339330467768
$ws.Range("S2").Value = 0.1938341213320902
$ws.Range("T2").Value = 0.1963983345368606

# Row 3
$ws.Range("M3").Value = 0.4493773333333334
$ws.Range("O3").Value = 0.1524833865370579
$ws.Range("P3").Value = 0.1545005747935906
$ws.Range("Q3").Value = 0.02966010233955556
$ws.Range("S3").Value = 0.1524833865370579
$ws.Range("T3").Value = 0.1545005747935906

# Row 4
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.185848666666667
$ws.Range("N4").Value = 3.557546
$ws.Range("O4").Value = 0.4023839370635547
$ws.Range("P4").Value = 0.4077070359984327
$ws.Range("Q4").Value = 0.0782691742631111
$ws.Range("R4").Value = 0.704422568368
$ws.Range("S4").Value = 0.4023839370635547
$ws.Range("T4").Value = 0.4077070359984327

# Row 5
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.5
$ws.Range("M5").Value = 0.115432
$ws.Range("N5").Value = 0.230864
$ws.Range("O5").Value = 0.03916855828972016
$ws.Range("P5").Value = 0.0264578102879744
$ws.Range("Q5").Value = 0.007618819818666667
$ws.Range("R5").Value = 0.045712918912
$ws.Range("S5").Value = 0.03916855828972016
$ws.Range("T5").Value = 0.0264578102879744

# Row 6
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.6251593333333333
$ws.Range("N6").Value = 1.875478
$ws.Range("O6").Value = 0.2121299967775769
$ws.Range("P6").Value = 0.2149362443831418
$ws.Range("Q6").Value = 0.04126218309155556
$ws.Range("R6").Value = 0.371359647824
$ws.Range("S6").Value = 0.2121299967775769
$ws.Range("T6").Value = 0.2149362443831418
